# Apply updates described by the diff:
# 1. Update the shared string for the Ponoko ordering note ($42.16 -> $46.53)
# 2. Update the unit cost in parts_1!C2 from 68 to 69.98 (this also recalculates
#    the SUM formula in C16)
# 3. Change the sheet selection on parts_1 to the full A1:D16 range

$wb = $excel.ActiveWorkbook

# --- Sheet "parts_1" edits ---
$ws = $wb.Worksheets.Item("parts_1")

# Update the unit cost; SUM(C2:C15) in C16 will recalc automatically.
$ws.Range("C2").Value = 69.98

# Update the Ponoko note text wherever it appears (it's a shared string, and
# the only cell using it is parts_10!G4, but search every sheet to be safe).
$oldNote = "Or order 10+ from Ponoko for @`$42.16"
$newNote = "Or order 10+ from Ponoko for @`$46.53"

foreach ($sheet in $wb.Worksheets) {
    $found = $sheet.Cells.Find($oldNote)
    while ($found -ne $null) {
        $found.Value = $newNote
        $found = $sheet.Cells.Find($oldNote)
    }
}

# Update the selection so the whole table A1:D16 is selected (no active-cell override)
$ws.Activate()
$ws.Range("A1:D16").Select()

$wb.Save()
